$d = $word.ActiveDocument

# 1) "On group page any user can see:" -> "On group page any visitor can see:"
$d.Content.Find.Execute(
    "On group page any user can see:", $true, $false, $false, $false, $false,
    $true, 1, $false, "On group page any visitor can see:", 2)

# 2) "Authorized user additional can see and use:" -> "User additionally can see and use:"
#    (only the part up to "use:" changes; stop short of the word "use:" itself so we
#    don't straddle the gramStart/gramEnd proofing-mark pair wrapped around it)
$d.Content.Find.Execute(
    "Authorized user additional can see and ", $true, $false, $false, $false, $false,
    $true, 1, $false, "User additionally can see and ", 2)

# 3) "Group owner additional can see and use: ... message." ->
#    "Group owner additionally can see and use: ... message. If Group owner leave group – group deletes. "
$d.Content.Find.Execute(
    "Group owner additional can see and use: settings button, button to exclude selected user, button to delete message on any posted message.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Group owner additionally can see and use: settings button, button to exclude selected user, button to delete message on any posted message. If Group owner leave group " + [char]0x2013 + " group deletes. ",
    2)

# 4) Move the "_GoBack" bookmark from the end of the document (after "PostgreSQL")
#    to right after the new sentence we just inserted, inside that same paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Content
$target.Find.Execute(
    "Group owner additionally can see and use: settings button, button to exclude selected user, button to delete message on any posted message. If Group owner leave group " + [char]0x2013 + " group deletes. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Word's collapsed-range-at-paragraph-end bookmark placement is unreliable, so
# temporarily insert a sentinel character to move the insertion point away
# from the exact paragraph-end boundary, add the bookmark, then remove the
# sentinel again. The bookmark stays correctly anchored in place.
$target.Collapse(0)
$target.InsertAfter("X")

$p = $target.Paragraphs(1)
$bookmarkPos = $p.Range.End - 2
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinel = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$sentinel.Delete()
